# with_cosine_similarities.xlsx
# - Fix header labels (they were off by one column relative to the data)
# - Fix a typo in H2 ("absord" -> "absorb")
# - Populate the new CosineSimilarity_1 column (J) with computed values
# - Shift the old CosineSimilarity_1/_2 data so CosineSimilarity_2 ends up in
#   column K, then drop the now-duplicate column L entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the header row (D1:K1); L1 will be removed by the column delete below ---
$ws.Range("D1").Value = "sentence"
$ws.Range("E1").Value = "verb_idx"
$ws.Range("F1").Value = "label"
$ws.Range("G1").Value = "word_1"
$ws.Range("H1").Value = "word_2"
$ws.Range("I1").Value = "word_3"
$ws.Range("J1").Value = "CosineSimilarity_1"
$ws.Range("K1").Value = "CosineSimilarity_2"

# --- Fix typo'd word_1 value in row 2 ---
$ws.Range("H2").Value = "absorb"

# --- New CosineSimilarity_1 values (column J), rows 2-30 ---
$j = @(0.2566680819897257, 0.2566680819897257, 0.3328949117932137, 0.4044853740991907, `
       0.2566680819897257, 0.2566680819897257, 0, 0, 0, 0.4652179322409485, `
       0.5271750183652546, 0.4385774836555155, 0.4393237764209459, 0.5605017618147151, `
       0.2428670676060369, 0.4393237764209459, 0.3408840893203698, 0.1903674402165263, `
       0.3722225428879729, 0.5664358014308725, 0.8241897019010316, 0.5089019823454131, `
       0.1023919432080509, 0.2169892474520203, 0.3949646915641254, 0.4925992310255374, `
       0.104213559963506, 0.4925992310255374, 0.3103177150743683)

# --- CosineSimilarity_2 values (column K), rows 2-30 ---
$k = @(0.3566897535651615, 0.5655341993882779, 0.3328949117932137, 0.2710425393070768, `
       0.4620813639364078, 0, 0, 0, 0, 0.5271750183652546, `
       0.2160821979520844, 0.3492648972773089, 0.5605017618147151, 0.3398623769483824, `
       0, 0.2068151976367784, 0.5283742105256388, 0.08354157338912775, `
       0.6955257109994322, 0, 0.7281898973172681, 0.5108451962195149, `
       0.4013323617022618, 0.1256868588701296, 0.2113993625174003, 0.51716012701194, `
       0, 0.4752230991641072, 0.1709367300599205)

for ($i = 0; $i -lt $j.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $j[$i]
    $ws.Cells.Item($row, 11).Value = $k[$i]
}

# --- Drop column L (its data has been folded into column K above) ---
$ws.Columns.Item(12).Delete()
